$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The experiment parameters were regenerated: distances D51/D64/D80 became
# D55/D69/D86 and the "S30" size level became "S31" everywhere it appears
# (condition labels, filenames, and the standalone lookup values).
$ws.Cells.Replace("D51", "D55", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$ws.Cells.Replace("D64", "D69", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$ws.Cells.Replace("D80", "D86", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
$ws.Cells.Replace("S30", "S31", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
